# Generate Report for Handoff
# Adds a new tracked file (f346e95e-57e6-4aca-bea4-ce0ef35182e2.md) as row 3
# on the Overview / zh-cn / de-de worksheets, mirroring the existing
# 495885ac-... row, and grows the three tables + autofilters to match.

$wb = $excel.ActiveWorkbook

$commitSha = "97a4161d0050ca3d66de2fd718d216758c79230f"
$newBase   = "f346e95e-57e6-4aca-bea4-ce0ef35182e2"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(3, 1).Value = "$newBase.md"
$wsOverview.Cells.Item(3, 2).Value = "e2e\$newBase.md"
$wsOverview.Cells.Item(3, 3).Value = ".md"
$wsOverview.Cells.Item(3, 4).Value = "'"
$wsOverview.Cells.Item(3, 4).Style = "Normal"
$wsOverview.Cells.Item(3, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 7).Value = "2016-08-26 08:43:57"
$wsOverview.Cells.Item(3, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$null = $wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$newBase.md",
    $null,
    $null,
    "e2e\$newBase.md"
)

$loOverview = $wsOverview.ListObjects.Item("Overview")
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Cells.Item(3, 1).Value = "$newBase.md"
$wsZhCn.Cells.Item(3, 2).Value = ".md"
$wsZhCn.Cells.Item(3, 3).Value = "Ready for handoff"
$wsZhCn.Cells.Item(3, 4).Value = "e2e"
$wsZhCn.Cells.Item(3, 5).Value = "ht"
$wsZhCn.Cells.Item(3, 6).Value = "'False"
$wsZhCn.Cells.Item(3, 6).Style = "Normal"
$wsZhCn.Cells.Item(3, 7).Value = "$newBase.0f17ffdc0225245bda8107a5962ed939b9b37863.zh-cn.xlf"
$wsZhCn.Cells.Item(3, 8).Value = "2016-08-26 08:43:53"
$wsZhCn.Cells.Item(3, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item(3, 9).Value = "'"
$wsZhCn.Cells.Item(3, 9).Style = "Normal"
$wsZhCn.Cells.Item(3, 10).Value = "'"
$wsZhCn.Cells.Item(3, 10).Style = "Normal"
$wsZhCn.Cells.Item(3, 11).Value = "0001-01-01 00:00:00"
$wsZhCn.Cells.Item(3, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item(3, 12).Value = "'"
$wsZhCn.Cells.Item(3, 12).Style = "Normal"
$wsZhCn.Cells.Item(3, 13).Value = "'True"
$wsZhCn.Cells.Item(3, 13).Style = "Normal"
$wsZhCn.Cells.Item(3, 14).Value = "'"
$wsZhCn.Cells.Item(3, 14).Style = "Normal"
$wsZhCn.Cells.Item(3, 15).Value = "'False"
$wsZhCn.Cells.Item(3, 15).Style = "Normal"
$wsZhCn.Cells.Item(3, 16).Value = "'"
$wsZhCn.Cells.Item(3, 16).Style = "Normal"

$null = $wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$newBase.md",
    $null,
    $null,
    "$newBase.md"
)

$loZhCn = $wsZhCn.ListObjects.Item("zh-cn")
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Cells.Item(3, 1).Value = "$newBase.md"
$wsDeDe.Cells.Item(3, 2).Value = ".md"
$wsDeDe.Cells.Item(3, 3).Value = "Ready for handoff"
$wsDeDe.Cells.Item(3, 4).Value = "e2e"
$wsDeDe.Cells.Item(3, 5).Value = "ht"
$wsDeDe.Cells.Item(3, 6).Value = "'False"
$wsDeDe.Cells.Item(3, 6).Style = "Normal"
$wsDeDe.Cells.Item(3, 7).Value = "$newBase.0f17ffdc0225245bda8107a5962ed939b9b37863.de-de.xlf"
$wsDeDe.Cells.Item(3, 8).Value = "2016-08-26 08:43:57"
$wsDeDe.Cells.Item(3, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item(3, 9).Value = "'"
$wsDeDe.Cells.Item(3, 9).Style = "Normal"
$wsDeDe.Cells.Item(3, 10).Value = "'"
$wsDeDe.Cells.Item(3, 10).Style = "Normal"
$wsDeDe.Cells.Item(3, 11).Value = "0001-01-01 00:00:00"
$wsDeDe.Cells.Item(3, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item(3, 12).Value = "'"
$wsDeDe.Cells.Item(3, 12).Style = "Normal"
$wsDeDe.Cells.Item(3, 13).Value = "'True"
$wsDeDe.Cells.Item(3, 13).Style = "Normal"
$wsDeDe.Cells.Item(3, 14).Value = "'"
$wsDeDe.Cells.Item(3, 14).Style = "Normal"
$wsDeDe.Cells.Item(3, 15).Value = "'False"
$wsDeDe.Cells.Item(3, 15).Style = "Normal"
$wsDeDe.Cells.Item(3, 16).Value = "'"
$wsDeDe.Cells.Item(3, 16).Style = "Normal"

$null = $wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$newBase.md",
    $null,
    $null,
    "$newBase.md"
)

$loDeDe = $wsDeDe.ListObjects.Item("de-de")
$loDeDe.Resize($wsDeDe.Range("A1:P3"))

Write-Host "Row 3 added to Overview, zh-cn and de-de sheets; tables resized."
